$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.155477046966553
$ws.Range("B1").Value = 3.395013093948364
$ws.Range("C1").Value = 3.10836124420166
$ws.Range("D1").Value = 3.529304504394531
$ws.Range("E1").Value = 1.562574028968811
